$d = $word.ActiveDocument

# "Elle sera constituée des élément" + (bookmark) + "s suivants, sous
# réserve de modifications si nécessaire :"
#   -> "Elle sera constituée des éléments suivants, sous réserve de
#       modifications si nécessaire" + " :" + (bookmark, now at the end)

# 1) Grow the first run's wording in place.
$d.Content.Find.Execute(
    "Elle sera constituée des élément",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Elle sera constituée des éléments suivants, sous réserve de modifications si nécessaire",
    2) | Out-Null

# 2) Shrink the old trailing run (which sat right after the hidden
#    _GoBack bookmark) down to just the punctuation that should remain.
$d.Content.Find.Execute(
    "s suivants, sous réserve de modifications si nécessaire :",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " :",
    2) | Out-Null

# Locate the paragraph that now holds the merged sentence.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Elle sera constituée*") {
        $target = $p
    }
}

# 3) Move the (hidden) _GoBack bookmark from between the two runs to the
#    very end of the paragraph (right before the paragraph mark), so it
#    ends up after both runs, matching the target layout.
$desiredPos = $target.Range.End - 1

# Moving a point bookmark to the position immediately before a paragraph
# mark is flaky when done directly, so first push that position one
# character away from the paragraph end by inserting a throwaway marker
# after it ...
$marker = $d.Range($desiredPos, $desiredPos)
$marker.InsertAfter("~")

# ... now retarget the bookmark at the (no longer boundary-adjacent)
# desired position ...
$d.Bookmarks.Add("_GoBack", $d.Range($desiredPos, $desiredPos)) | Out-Null
$bm = $d.Bookmarks.Item("_GoBack")

# ... and finally remove the throwaway marker that now trails it.
$trailing = $d.Range($bm.End, $bm.End + 1)
$trailing.Delete()
